$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "Day 5" (第五天) column of the third weekly block (rows 44-57)
# with the same kind of sleep-diary answers used elsewhere in the sheet.
$ws.Range("E44").Value = "8：30"
$ws.Range("E45").Value = "8：45"
$ws.Range("E46").Value = "23：10"
$ws.Range("E47").Value = "23：20"
$ws.Range("E48").Value = 10
$ws.Range("E49").Value = 2
$ws.Range("E50").Value = 80
$ws.Range("E51").Value = 460
$ws.Range("E52").Value = "无"
$ws.Range("E53").Value = 20
$ws.Range("E54").Value = 2
$ws.Range("E55").Value = 4
$ws.Range("E56").Value = 2
$ws.Range("E57").Value = "无"

# Leave the cursor where the author ended up after typing the last entry.
$ws.Range("E57").Select() | Out-Null
